# Updated the validationMain and the way analytic results are obtained.
# Refresh the node/date data pulled from Postgres (replacing the old IHPC
# query results) on the "nodes for testing auto-validation" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (nid, date) as pulled from Postgres - rows 2..61
$data = @(
    @(509057, 20150930),
    @(509057, 20151001),
    @(509059, 20150930),
    @(509059, 20151001),
    @(512116, 20150929),
    @(512117, 20150930),
    @(512118, 20150930),
    @(512118, 20151001),
    @(512119, 20150929),
    @(512119, 20150930),
    @(512123, 20150930),
    @(512123, 20151001),
    @(512124, 20150929),
    @(512124, 20150930),
    @(512124, 20151001),
    @(513425, 20150929),
    @(513425, 20150930),
    @(513425, 20151001),
    @(513427, 20150930),
    @(513432, 20150930),
    @(513435, 20150929),
    @(513435, 20150930),
    @(513435, 20151001),
    @(513437, 20150929),
    @(513471, 20150929),
    @(513471, 20150930),
    @(513471, 20151001),
    @(513472, 20150930),
    @(513473, 20150929),
    @(513473, 20150930),
    @(513475, 20150929),
    @(513475, 20150930),
    @(513475, 20151001),
    @(513476, 20150929),
    @(513476, 20150930),
    @(513478, 20150929),
    @(514500, 20150929),
    @(514500, 20150930),
    @(514501, 20150929),
    @(514503, 20150929),
    @(514503, 20150930),
    @(514505, 20150929),
    @(514505, 20150930),
    @(514511, 20151001),
    @(514511, 20151002),
    @(514594, 20150929),
    @(514594, 20150930),
    @(514594, 20151001),
    @(514594, 20151002),
    @(514595, 20150929),
    @(514595, 20151002),
    @(514596, 20150929),
    @(514598, 20150930),
    @(517606, 20150929),
    @(517610, 20150929),
    @(517612, 20150929),
    @(517612, 20151002),
    @(517617, 20150929),
    @(517617, 20150930),
    @(517620, 20150930)
)

$oldLastRow = 64
$newLastRow = $data.Count + 1

# Write the refreshed data into rows 2..(newLastRow)
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# The new data set is shorter than before, so clear out the now-unused
# trailing rows that used to hold values (rows 62..64).
if ($newLastRow -lt $oldLastRow) {
    $clearRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($oldLastRow, 2))
    $clearRange.ClearContents()
}

# Update the active selection to match the saved view state.
$ws.Range("J9").Select()
